# chore: update Sheets via scheduled runner
# Refreshes cached marketboard/profit figures (currentAveragePrice*, LevePrice*, LeveProfit*)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR profit sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 194.35294
$ws.Range("I33").Value = 197.125
$ws.Range("K33").Value = 197.125
$ws.Range("M33").Value = 31.875
$ws.Range("H43").Value = 2937.7778
$ws.Range("J43").Value = 2168
$ws.Range("L43").Value = 2168
$ws.Range("N43").Value = -2306
$ws.Range("H51").Value = 9270
$ws.Range("J51").Value = 8691.4
$ws.Range("L51").Value = 8691.4
$ws.Range("N51").Value = -9659.4
$ws.Range("H58").Value = 812.0909
$ws.Range("J58").Value = 2419
$ws.Range("L58").Value = 7257
$ws.Range("N58").Value = -7557
$ws.Range("H98").Value = 1851.1915
$ws.Range("I98").Value = 1853.4889
$ws.Range("K98").Value = 1853.4889
$ws.Range("M98").Value = -355.4889000000001
$ws.Range("H122").Value = 1851.1915
$ws.Range("I122").Value = 1853.4889
$ws.Range("K122").Value = 5560.4667
$ws.Range("M122").Value = -3110.4667
$ws.Range("H132").Value = 7165.2173
$ws.Range("I132").Value = 7165.2173
$ws.Range("K132").Value = 21495.6519
$ws.Range("M132").Value = -18965.6519
$ws.Range("H138").Value = 350565.22
$ws.Range("J138").Value = 523890.84
$ws.Range("L138").Value = 1571672.52
$ws.Range("N138").Value = -1581952.52
$ws.Range("H141").Value = 5986.696
$ws.Range("I141").Value = 3699.7144
$ws.Range("K141").Value = 11099.1432
$ws.Range("M141").Value = -5919.143199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 851.4706
$ws.Range("I2").Value = 716.1429000000001
$ws.Range("K2").Value = 716.1429000000001
$ws.Range("M2").Value = -603.1429000000001
$ws.Range("H74").Value = 215399.39
$ws.Range("I74").Value = 309886.72
$ws.Range("K74").Value = 309886.72
$ws.Range("M74").Value = -309012.72
$ws.Range("H77").Value = 215399.39
$ws.Range("I77").Value = 309886.72
$ws.Range("K77").Value = 1549433.6
$ws.Range("M77").Value = -1545065.6
$ws.Range("H110").Value = 2965.625
$ws.Range("I110").Value = 1771.2
$ws.Range("K110").Value = 1771.2
$ws.Range("M110").Value = 273.8
$ws.Range("H116").Value = 851.4706
$ws.Range("I116").Value = 716.1429000000001
$ws.Range("K116").Value = 716.1429000000001
$ws.Range("M116").Value = 1577.8571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 851.4706
$ws.Range("I3").Value = 716.1429000000001
$ws.Range("K3").Value = 716.1429000000001
$ws.Range("M3").Value = -602.1429000000001
$ws.Range("H64").Value = 695.625
$ws.Range("J64").Value = 920
$ws.Range("L64").Value = 920
$ws.Range("N64").Value = -1370
$ws.Range("H67").Value = 695.625
$ws.Range("J67").Value = 920
$ws.Range("L67").Value = 920
$ws.Range("N67").Value = -2480
$ws.Range("H134").Value = 4583.1665
$ws.Range("I134").Value = 3856.8572
$ws.Range("K134").Value = 11570.5716
$ws.Range("M134").Value = -9035.571599999999
$ws.Range("H141").Value = 71495.60000000001
$ws.Range("I141").Value = 40000
$ws.Range("K141").Value = 40000
$ws.Range("M141").Value = -34820

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5744.273
$ws.Range("I31").Value = 4350.643
$ws.Range("K31").Value = 4350.643
$ws.Range("M31").Value = -4055.643
$ws.Range("H34").Value = 5744.273
$ws.Range("I34").Value = 4350.643
$ws.Range("K34").Value = 4350.643
$ws.Range("M34").Value = -4148.643
$ws.Range("H58").Value = 2319.9487
$ws.Range("I58").Value = 1529
$ws.Range("J58").Value = 3732.3572
$ws.Range("K58").Value = 1529
$ws.Range("L58").Value = 3732.3572
$ws.Range("M58").Value = -1326
$ws.Range("N58").Value = -4138.3572
$ws.Range("H134").Value = 3282.9697
$ws.Range("I134").Value = 3043.1614
$ws.Range("K134").Value = 9129.484199999999
$ws.Range("M134").Value = -6594.484199999999
$ws.Range("H136").Value = 2319.9487
$ws.Range("I136").Value = 1529
$ws.Range("J136").Value = 3732.3572
$ws.Range("K136").Value = 4587
$ws.Range("L136").Value = 11197.0716
$ws.Range("M136").Value = -2037
$ws.Range("N136").Value = -16297.0716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 3381.6365
$ws.Range("J64").Value = 3529.9
$ws.Range("L64").Value = 10589.7
$ws.Range("N64").Value = -11129.7
$ws.Range("H67").Value = 3381.6365
$ws.Range("J67").Value = 3529.9
$ws.Range("L67").Value = 10589.7
$ws.Range("N67").Value = -12461.7
$ws.Range("H70").Value = 4337.3335
$ws.Range("I70").Value = 498
$ws.Range("K70").Value = 1494
$ws.Range("M70").Value = -1179
$ws.Range("H73").Value = 4337.3335
$ws.Range("I73").Value = 498
$ws.Range("K73").Value = 1494
$ws.Range("M73").Value = -402
$ws.Range("H101").Value = 4874.75
$ws.Range("I101").Value = 1499.5
$ws.Range("J101").Value = 8250
$ws.Range("K101").Value = 4498.5
$ws.Range("L101").Value = 24750
$ws.Range("M101").Value = -2064.5
$ws.Range("N101").Value = -29618
$ws.Range("H121").Value = 125484.75
$ws.Range("I121").Value = 155.8
$ws.Range("K121").Value = 467.4
$ws.Range("M121").Value = 842.5999999999999
$ws.Range("H122").Value = 1419.7894
$ws.Range("I122").Value = 918
$ws.Range("J122").Value = 1599
$ws.Range("K122").Value = 8262
$ws.Range("L122").Value = 14391
$ws.Range("M122").Value = -5812
$ws.Range("N122").Value = -19291
$ws.Range("H139").Value = 1662.25
$ws.Range("I139").Value = 1234.7
$ws.Range("K139").Value = 3704.1
$ws.Range("M139").Value = 1435.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 4104
$ws.Range("J36").Value = 4104
$ws.Range("L36").Value = 4104
$ws.Range("N36").Value = -5074
$ws.Range("H126").Value = 4541.9375
$ws.Range("J126").Value = 9614
$ws.Range("L126").Value = 28842
$ws.Range("N126").Value = -33782
$ws.Range("H132").Value = 8587.6
$ws.Range("I132").Value = 176.6
$ws.Range("J132").Value = 16998.6
$ws.Range("K132").Value = 529.8
$ws.Range("L132").Value = 50995.8
$ws.Range("M132").Value = 2000.2
$ws.Range("N132").Value = -56055.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3371.8
$ws.Range("I7").Value = 3335.389
$ws.Range("J7").Value = 3699.5
$ws.Range("K7").Value = 3335.389
$ws.Range("L7").Value = 3699.5
$ws.Range("M7").Value = -3223.389
$ws.Range("N7").Value = -3923.5
$ws.Range("H22").Value = 798
$ws.Range("I22").Value = 581.6667
$ws.Range("J22").Value = 983.4286
$ws.Range("K22").Value = 581.6667
$ws.Range("L22").Value = 983.4286
$ws.Range("M22").Value = -286.6667
$ws.Range("N22").Value = -1573.4286
$ws.Range("H27").Value = 798
$ws.Range("I27").Value = 581.6667
$ws.Range("J27").Value = 983.4286
$ws.Range("K27").Value = 581.6667
$ws.Range("L27").Value = 983.4286
$ws.Range("M27").Value = -474.6667
$ws.Range("N27").Value = -1197.4286
$ws.Range("H40").Value = 4950.4614
$ws.Range("I40").Value = 5008.7188
$ws.Range("K40").Value = 5008.7188
$ws.Range("M40").Value = -4872.7188
$ws.Range("H126").Value = 3371.8
$ws.Range("I126").Value = 3335.389
$ws.Range("J126").Value = 3699.5
$ws.Range("K126").Value = 10006.167
$ws.Range("L126").Value = 11098.5
$ws.Range("M126").Value = -7536.167000000001
$ws.Range("N126").Value = -16038.5
$ws.Range("H136").Value = 2843.804
$ws.Range("I136").Value = 2480.8948
$ws.Range("K136").Value = 7442.6844
$ws.Range("M136").Value = -4892.6844

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2178.9
$ws.Range("I126").Value = 2123.75
$ws.Range("K126").Value = 6371.25
$ws.Range("M126").Value = -3901.25
$ws.Range("H132").Value = 2337.1428
$ws.Range("I132").Value = 2203.1562
$ws.Range("K132").Value = 6609.4686
$ws.Range("M132").Value = -4079.4686
$ws.Range("H136").Value = 55559704
$ws.Range("I136").Value = 76924250
$ws.Range("K136").Value = 230772750
$ws.Range("M136").Value = -230770200
$ws.Range("H137").Value = 77499
$ws.Range("J137").Value = 77499
$ws.Range("L137").Value = 77499
$ws.Range("N137").Value = -87699
